# "update ismart for grant extension"
#
# The "grants" worksheet has a row for the "Paths to Success..." grant
# (status = Under Review) that needs to be removed entirely, and the
# I-SMART grant's status needs to move from "Previously Funded" to
# "Currently Funded" (the grant was extended).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grants")

# Remove the entire row for the "Paths to Success..." grant (row 3).
# This shifts the I-SMART row (was row 4) up to row 3, and the
# "Unfunded" row (was row 5) up to row 4.
$ws.Rows.Item(3).Delete()

# The I-SMART grant (now in row 3) was extended, so update its status
# from "Previously Funded" to "Currently Funded".
$ws.Range("A3").Value = "Currently Funded"

# Keep the active selection on the now-shorter data range.
$ws.Range("A4").Select()
